$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 55912.723
$ws.Range("I53").Value = 111271.445
$ws.Range("J53").Value = 554
$ws.Range("K53").Value = 111271.445
$ws.Range("L53").Value = 554
$ws.Range("M53").Value = -110634.445
$ws.Range("N53").Value = -1828

$ws.Range("H64").Value = 4215.067
$ws.Range("I64").Value = 4608.6
$ws.Range("J64").Value = 3428
$ws.Range("K64").Value = 4608.6
$ws.Range("L64").Value = 3428
$ws.Range("M64").Value = -4360.6
$ws.Range("N64").Value = -3924

$ws.Range("H67").Value = 4215.067
$ws.Range("I67").Value = 4608.6
$ws.Range("J67").Value = 3428
$ws.Range("K67").Value = 4608.6
$ws.Range("L67").Value = 3428
$ws.Range("M67").Value = -3750.6
$ws.Range("N67").Value = -5144

$ws.Range("H76").Value = 7006.1143
$ws.Range("I76").Value = 12473.272
$ws.Range("J76").Value = 4500.3335
$ws.Range("K76").Value = 12473.272
$ws.Range("L76").Value = 4500.3335
$ws.Range("M76").Value = -12158.272
$ws.Range("N76").Value = -5130.3335

$ws.Range("H79").Value = 7006.1143
$ws.Range("I79").Value = 12473.272
$ws.Range("J79").Value = 4500.3335
$ws.Range("K79").Value = 12473.272
$ws.Range("L79").Value = 4500.3335
$ws.Range("M79").Value = -11381.272
$ws.Range("N79").Value = -6684.3335

$ws.Range("H98").Value = 1479.8334
$ws.Range("I98").Value = 1311.3636
$ws.Range("K98").Value = 1311.3636
$ws.Range("M98").Value = 186.6364000000001

$ws.Range("H122").Value = 1479.8334
$ws.Range("I122").Value = 1311.3636
$ws.Range("K122").Value = 3934.0908
$ws.Range("M122").Value = -1484.0908

$ws.Range("H132").Value = 2022.6945
$ws.Range("I132").Value = 1779.0312
$ws.Range("J132").Value = 3972
$ws.Range("K132").Value = 5337.0936
$ws.Range("L132").Value = 11916
$ws.Range("M132").Value = -2807.0936
$ws.Range("N132").Value = -16976

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H74").Value = 1587
$ws.Range("I74").Value = 1262.4828
$ws.Range("J74").Value = 2310.923
$ws.Range("K74").Value = 1262.4828
$ws.Range("L74").Value = 2310.923
$ws.Range("M74").Value = -388.4828
$ws.Range("N74").Value = -4058.923

$ws.Range("H77").Value = 1587
$ws.Range("I77").Value = 1262.4828
$ws.Range("J77").Value = 2310.923
$ws.Range("K77").Value = 6312.414
$ws.Range("L77").Value = 11554.615
$ws.Range("M77").Value = -1944.414
$ws.Range("N77").Value = -20290.615

$ws.Range("H88").Value = 2515.2856
$ws.Range("J88").Value = 2681.4
$ws.Range("L88").Value = 2681.4
$ws.Range("N88").Value = -3493.4

$ws.Range("H91").Value = 2515.2856
$ws.Range("J91").Value = 2681.4
$ws.Range("L91").Value = 2681.4
$ws.Range("N91").Value = -5489.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1459.8889
$ws.Range("I5").Value = 1019.8571
$ws.Range("K5").Value = 1019.8571
$ws.Range("M5").Value = -906.8570999999999

$ws.Range("H105").Value = 5700
$ws.Range("I105").Value = 7260
$ws.Range("J105").Value = 4400
$ws.Range("K105").Value = 7260
$ws.Range("L105").Value = 4400
$ws.Range("M105").Value = -5513
$ws.Range("N105").Value = -7894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2374.8235
$ws.Range("I16").Value = 2230
$ws.Range("J16").Value = 2503.5557
$ws.Range("K16").Value = 2230
$ws.Range("L16").Value = 2503.5557
$ws.Range("M16").Value = -1943
$ws.Range("N16").Value = -3077.5557

$ws.Range("H62").Value = 6472.778
$ws.Range("I62").Value = 6773.091
$ws.Range("J62").Value = 6000.857
$ws.Range("K62").Value = 6773.091
$ws.Range("L62").Value = 6000.857
$ws.Range("M62").Value = -6149.091
$ws.Range("N62").Value = -7248.857

$ws.Range("H65").Value = 6472.778
$ws.Range("I65").Value = 6773.091
$ws.Range("J65").Value = 6000.857
$ws.Range("K65").Value = 33865.455
$ws.Range("L65").Value = 30004.285
$ws.Range("M65").Value = -30745.455
$ws.Range("N65").Value = -36244.285

$ws.Range("H113").Value = 2374.8235
$ws.Range("I113").Value = 2230
$ws.Range("J113").Value = 2503.5557
$ws.Range("K113").Value = 2230
$ws.Range("L113").Value = 2503.5557
$ws.Range("M113").Value = -60
$ws.Range("N113").Value = -6843.5557

$ws.Range("H132").Value = 1590.2433
$ws.Range("I132").Value = 818.5
$ws.Range("J132").Value = 3015
$ws.Range("K132").Value = 2455.5
$ws.Range("L132").Value = 9045
$ws.Range("M132").Value = 74.5
$ws.Range("N132").Value = -14105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5297
$ws.Range("I70").Value = 5440.2
$ws.Range("J70").Value = 4581
$ws.Range("K70").Value = 5440.2
$ws.Range("L70").Value = 4581
$ws.Range("M70").Value = -5170.2
$ws.Range("N70").Value = -5121

$ws.Range("H73").Value = 5297
$ws.Range("I73").Value = 5440.2
$ws.Range("J73").Value = 4581
$ws.Range("K73").Value = 5440.2
$ws.Range("L73").Value = 4581
$ws.Range("M73").Value = -4504.2
$ws.Range("N73").Value = -6453

$ws.Range("H80").Value = 6367.36
$ws.Range("I80").Value = 8979
$ws.Range("K80").Value = 8979
$ws.Range("M80").Value = -7981

$ws.Range("H83").Value = 6367.36
$ws.Range("I83").Value = 8979
$ws.Range("K83").Value = 44895
$ws.Range("M83").Value = -39903

$ws.Range("H113").Value = 35715800
$ws.Range("I113").Value = 50001108
$ws.Range("J113").Value = 2526.625
$ws.Range("K113").Value = 50001108
$ws.Range("L113").Value = 2526.625
$ws.Range("M113").Value = -49998938
$ws.Range("N113").Value = -6866.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 14444
$ws.Range("J76").Value = 14444
$ws.Range("L76").Value = 14444
$ws.Range("N76").Value = -15120

$ws.Range("H79").Value = 14444
$ws.Range("J79").Value = 14444
$ws.Range("L79").Value = 14444
$ws.Range("N79").Value = -16784

$ws.Range("H114").Value = 23833.334
$ws.Range("J114").Value = 23833.334
$ws.Range("L114").Value = 23833.334
$ws.Range("N114").Value = -32511.334

$ws.Range("H133").Value = 170326
$ws.Range("J133").Value = 170326
$ws.Range("L133").Value = 170326
$ws.Range("N133").Value = -175386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 30173
$ws.Range("J76").Value = 30173
$ws.Range("L76").Value = 30173
$ws.Range("N76").Value = -30803

$ws.Range("H79").Value = 30173
$ws.Range("J79").Value = 30173
$ws.Range("L79").Value = 30173
$ws.Range("N79").Value = -32357

$ws.Range("H104").Value = 31000
$ws.Range("J104").Value = 31000
$ws.Range("L104").Value = 31000
$ws.Range("N104").Value = -37988

$ws.Range("H122").Value = 2179
$ws.Range("I122").Value = 1589.25
$ws.Range("J122").Value = 3751.6667
$ws.Range("K122").Value = 4767.75
$ws.Range("L122").Value = 11255.0001
$ws.Range("M122").Value = -2317.75
$ws.Range("N122").Value = -16155.0001

$ws.Range("H126").Value = 938.3889
$ws.Range("I126").Value = 779.1429000000001
$ws.Range("J126").Value = 1495.75
$ws.Range("K126").Value = 2337.4287
$ws.Range("L126").Value = 4487.25
$ws.Range("M126").Value = 132.5712999999996
$ws.Range("N126").Value = -9427.25
